$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").NumberFormat = "@"
$ws.Range("A8").Value = 'Volume 31   Number  20'
$ws.Range("A8").NumberFormat = 'General'

$ws.Range("C9").NumberFormat = "@"
$ws.Range("C9").Value = 'Report Covering the Week  5/13/2024  Through  5/19/2024'
$ws.Range("C9").NumberFormat = 'General'

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0'
$ws.Range("D14").NumberFormat = 'General'

$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '***.*'
$ws.Range("E14").NumberFormat = 'General'

$ws.Range("N14").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("N14").Value = -84.615384615384

$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = '0'
$ws.Range("C15").NumberFormat = 'General'

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0'
$ws.Range("D15").NumberFormat = 'General'

$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '***.*'
$ws.Range("E15").NumberFormat = 'General'

$ws.Range("F15").NumberFormat = '#,##0'
$ws.Range("F15").Value = 1

$ws.Range("H15").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("H15").Value = -66.666666666666

$ws.Range("I15").NumberFormat = '#,##0'
$ws.Range("I15").Value = 5

$ws.Range("K15").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("K15").Value = -44.444444444444

$ws.Range("L15").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("L15").Value = -37.5

$ws.Range("M15").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("M15").Value = -54.545454545454

$ws.Range("N15").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("N15").Value = -54.545454545454

$ws.Range("C16").NumberFormat = '#,##0'
$ws.Range("C16").Value = 3

$ws.Range("D16").NumberFormat = '#,##0'
$ws.Range("D16").Value = 4

$ws.Range("E16").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("E16").Value = -25

$ws.Range("F16").NumberFormat = '#,##0'
$ws.Range("F16").Value = 19

$ws.Range("H16").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("H16").Value = 58.333333333333

$ws.Range("I16").NumberFormat = '#,##0'
$ws.Range("I16").Value = 79

$ws.Range("J16").NumberFormat = '#,##0'
$ws.Range("J16").Value = 76

$ws.Range("K16").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("K16").Value = 3.947368421052

$ws.Range("L16").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("L16").Value = 3.947368421052

$ws.Range("M16").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("M16").Value = -12.222222222222

$ws.Range("N16").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("N16").Value = -76.488095238095

$ws.Range("C17").NumberFormat = '#,##0'
$ws.Range("C17").Value = 6

$ws.Range("D17").NumberFormat = '#,##0'
$ws.Range("D17").Value = 7

$ws.Range("E17").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("E17").Value = -14.285714285714

$ws.Range("F17").NumberFormat = '#,##0'
$ws.Range("F17").Value = 34

$ws.Range("G17").NumberFormat = '#,##0'
$ws.Range("G17").Value = 20

$ws.Range("H17").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("H17").Value = 70

$ws.Range("I17").NumberFormat = '#,##0'
$ws.Range("I17").Value = 137

$ws.Range("J17").NumberFormat = '#,##0'
$ws.Range("J17").Value = 116

$ws.Range("K17").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("K17").Value = 18.103448275862

$ws.Range("L17").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("L17").Value = 28.03738317757

$ws.Range("M17").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("M17").Value = 153.703703703704

$ws.Range("N17").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("N17").Value = -0.724637681159

$ws.Range("C18").NumberFormat = '#,##0'
$ws.Range("C18").Value = 1

$ws.Range("E18").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("E18").Value = 0

$ws.Range("F18").NumberFormat = '#,##0'
$ws.Range("F18").Value = 8

$ws.Range("G18").NumberFormat = '#,##0'
$ws.Range("G18").Value = 6

$ws.Range("H18").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("H18").Value = 33.333333333333

$ws.Range("I18").NumberFormat = '#,##0'
$ws.Range("I18").Value = 47

$ws.Range("J18").NumberFormat = '#,##0'
$ws.Range("J18").Value = 38

$ws.Range("K18").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("K18").Value = 23.684210526315

$ws.Range("L18").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("L18").Value = -6

$ws.Range("M18").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("M18").Value = -51.041666666666

$ws.Range("N18").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("N18").Value = -89.220183486238

$ws.Range("C19").NumberFormat = '#,##0'
$ws.Range("C19").Value = 8

$ws.Range("E19").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("E19").Value = -33.333333333333

$ws.Range("F19").NumberFormat = '#,##0'
$ws.Range("F19").Value = 35

$ws.Range("G19").NumberFormat = '#,##0'
$ws.Range("G19").Value = 55

$ws.Range("H19").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("H19").Value = -36.363636363636

$ws.Range("I19").NumberFormat = '#,##0'
$ws.Range("I19").Value = 195

$ws.Range("J19").NumberFormat = '#,##0'
$ws.Range("J19").Value = 235

$ws.Range("K19").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("K19").Value = -17.021276595744

$ws.Range("L19").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("L19").Value = -29.090909090909

$ws.Range("M19").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("M19").Value = 58.536585365853

$ws.Range("N19").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("N19").Value = -9.302325581395

$ws.Range("C20").NumberFormat = '#,##0'
$ws.Range("C20").Value = 6

$ws.Range("D20").NumberFormat = '#,##0'
$ws.Range("D20").Value = 4

$ws.Range("E20").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("E20").Value = 50

$ws.Range("F20").NumberFormat = '#,##0'
$ws.Range("F20").Value = 24

$ws.Range("G20").NumberFormat = '#,##0'
$ws.Range("G20").Value = 10

$ws.Range("H20").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("H20").Value = 140

$ws.Range("I20").NumberFormat = '#,##0'
$ws.Range("I20").Value = 122

$ws.Range("J20").NumberFormat = '#,##0'
$ws.Range("J20").Value = 81

$ws.Range("K20").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("K20").Value = 50.617283950617

$ws.Range("L20").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("L20").Value = 19.607843137254

$ws.Range("M20").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("M20").Value = 7.964601769911

$ws.Range("N20").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("N20").Value = -90.665646518745

$ws.Range("C21").NumberFormat = '#,##0'
$ws.Range("C21").Value = 24

$ws.Range("D21").NumberFormat = '#,##0'
$ws.Range("D21").Value = 28

$ws.Range("E21").NumberFormat = '#,##0.00;"-"#,##0.00'
$ws.Range("E21").Value = -14.285714285714

$ws.Range("F21").NumberFormat = '#,##0'
$ws.Range("F21").Value = 122

$ws.Range("G21").NumberFormat = '#,##0'
$ws.Range("G21").Value = 107

$ws.Range("H21").NumberFormat = '#,##0.00;"-"#,##0.00'
$ws.Range("H21").Value = 14.018691588785

$ws.Range("I21").NumberFormat = '#,##0'
$ws.Range("I21").Value = 587

$ws.Range("J21").NumberFormat = '#,##0'
$ws.Range("J21").Value = 556

$ws.Range("K21").NumberFormat = '#,##0.00;"-"#,##0.00'
$ws.Range("K21").Value = 5.575539568345

$ws.Range("L21").NumberFormat = '#,##0.00;"-"#,##0.00'
$ws.Range("L21").Value = -5.169628432956

$ws.Range("M21").NumberFormat = '#,##0.00;"-"#,##0.00'
$ws.Range("M21").Value = 20.286885245901

$ws.Range("N21").NumberFormat = '#,##0.00;"-"#,##0.00'
$ws.Range("N21").Value = -76.099348534202

$ws.Range("C22").NumberFormat = '#,##0'
$ws.Range("C22").Value = 3

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0'
$ws.Range("D22").NumberFormat = 'General'

$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '***.*'
$ws.Range("E22").NumberFormat = 'General'

$ws.Range("F22").NumberFormat = '#,##0'
$ws.Range("F22").Value = 3

$ws.Range("H22").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("H22").Value = 0

$ws.Range("I22").NumberFormat = '#,##0'
$ws.Range("I22").Value = 8

$ws.Range("K22").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("K22").Value = 0

$ws.Range("L22").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("L22").Value = 60

$ws.Range("M22").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("M22").Value = -20

$ws.Range("D24").NumberFormat = '#,##0'
$ws.Range("D24").Value = 33

$ws.Range("E24").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("E24").Value = -39.393939393939

$ws.Range("F24").NumberFormat = '#,##0'
$ws.Range("F24").Value = 87

$ws.Range("G24").NumberFormat = '#,##0'
$ws.Range("G24").Value = 95

$ws.Range("H24").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("H24").Value = -8.421052631578

$ws.Range("I24").NumberFormat = '#,##0'
$ws.Range("I24").Value = 463

$ws.Range("J24").NumberFormat = '#,##0'
$ws.Range("J24").Value = 481

$ws.Range("K24").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("K24").Value = -3.742203742203

$ws.Range("L24").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("L24").Value = -15.201465201465

$ws.Range("M24").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("M24").Value = 91.322314049586

$ws.Range("C25").NumberFormat = '#,##0'
$ws.Range("C25").Value = 8

$ws.Range("D25").NumberFormat = '#,##0'
$ws.Range("D25").Value = 13

$ws.Range("E25").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("E25").Value = -38.461538461538

$ws.Range("F25").NumberFormat = '#,##0'
$ws.Range("F25").Value = 49

$ws.Range("G25").NumberFormat = '#,##0'
$ws.Range("G25").Value = 30

$ws.Range("H25").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("H25").Value = 63.333333333333

$ws.Range("I25").NumberFormat = '#,##0'
$ws.Range("I25").Value = 229

$ws.Range("J25").NumberFormat = '#,##0'
$ws.Range("J25").Value = 196

$ws.Range("K25").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("K25").Value = 16.836734693877

$ws.Range("L25").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("L25").Value = -3.37552742616

$ws.Range("C26").NumberFormat = '#,##0'
$ws.Range("C26").Value = 13

$ws.Range("D26").NumberFormat = '#,##0'
$ws.Range("D26").Value = 13

$ws.Range("E26").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("E26").Value = 0

$ws.Range("F26").NumberFormat = '#,##0'
$ws.Range("F26").Value = 54

$ws.Range("G26").NumberFormat = '#,##0'
$ws.Range("G26").Value = 50

$ws.Range("H26").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("H26").Value = 8

$ws.Range("I26").NumberFormat = '#,##0'
$ws.Range("I26").Value = 210

$ws.Range("J26").NumberFormat = '#,##0'
$ws.Range("J26").Value = 203

$ws.Range("K26").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("K26").Value = 3.448275862068

$ws.Range("L26").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("L26").Value = 28.048780487804

$ws.Range("M26").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("M26").Value = 14.130434782608

$ws.Range("C27").NumberFormat = "@"
$ws.Range("C27").Value = '0'
$ws.Range("C27").NumberFormat = 'General'

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0'
$ws.Range("D27").NumberFormat = 'General'

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '***.*'
$ws.Range("E27").NumberFormat = 'General'

$ws.Range("F27").NumberFormat = '#,##0'
$ws.Range("F27").Value = 2

$ws.Range("H27").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("H27").Value = -33.333333333333

$ws.Range("I27").NumberFormat = '#,##0'
$ws.Range("I27").Value = 9

$ws.Range("K27").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("K27").Value = -30.76923076923

$ws.Range("L27").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("L27").Value = -35.714285714285

$ws.Range("C28").NumberFormat = '#,##0'
$ws.Range("C28").Value = 1

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0'
$ws.Range("D28").NumberFormat = 'General'

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '***.*'
$ws.Range("E28").NumberFormat = 'General'

$ws.Range("G28").NumberFormat = '#,##0'
$ws.Range("G28").Value = 2

$ws.Range("H28").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("H28").Value = 200

$ws.Range("I28").NumberFormat = '#,##0'
$ws.Range("I28").Value = 17

$ws.Range("K28").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("K28").Value = -22.727272727272

$ws.Range("L28").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("L28").Value = -29.166666666666

$ws.Range("N29").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("N29").Value = -95

$ws.Range("N30").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("N30").Value = -95

$ws.Range("D33").NumberFormat = '#,##0'
$ws.Range("D33").Value = 1

$ws.Range("E33").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("E33").Value = -100

$ws.Range("G33").NumberFormat = '#,##0'
$ws.Range("G33").Value = 1

$ws.Range("H33").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("H33").Value = -100

$ws.Range("J33").NumberFormat = '#,##0'
$ws.Range("J33").Value = 3

$ws.Range("K33").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("K33").Value = -33.333333333333
